# Adds 20 new response rows (82-101) to the "historico" sheet, matching
# the rows logged between 29/12/2025 21:15 and 29/12/2025 22:11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, A(id), B(data/timestamp), C(questao_id), C is text?, D(disciplina), E(assunto), F(resultado)
$rows = @(
    @(82, 81, "29/12/2025 21:15", "955", $true,  "Inglês", "Verbs", 1),
    @(83, 82, "29/12/2025 21:16", "954", $true,  "Inglês", "Pronouns", 1),
    @(84, 83, "29/12/2025 21:16", "953", $true,  "Inglês", "Semantic", 1),
    @(85, 84, "29/12/2025 21:37", 952,   $false, "Inglês", "Semantic", 1),
    @(86, 85, "29/12/2025 21:38", 951,   $false, "Inglês", "Interpretação de Texto", 1),
    @(87, 86, "29/12/2025 21:40", "980", $true,  "Inglês", "Interpretação de Texto", 1),
    @(88, 87, "29/12/2025 21:41", 979,   $false, "Inglês", "Interpretação de Texto", 0),
    @(89, 88, "29/12/2025 21:43", 978,   $false, "Inglês", "Interpretação de Texto", 1),
    @(90, 89, "29/12/2025 21:44", 977,   $false, "Inglês", "Semantic", 0),
    @(91, 90, "29/12/2025 21:47", 976,   $false, "Inglês", "Semantic", 1),
    @(92, 91, "29/12/2025 21:47", "975", $true,  "Inglês", "Pronouns", 0),
    @(93, 92, "29/12/2025 21:49", 974,   $false, "Inglês", "Interpretação de Texto", 1),
    @(94, 93, "29/12/2025 21:50", 973,   $false, "Inglês", "Semantic", 1),
    @(95, 94, "29/12/2025 21:52", 972,   $false, "Inglês", "Semantic", 1),
    @(96, 95, "29/12/2025 21:53", 971,   $false, "Inglês", "Semantic", 1),
    @(97, 96, "29/12/2025 21:56", "881", $true,  "Inglês", "Semantic", 1),
    @(98, 97, "29/12/2025 21:57", "880", $true,  "Inglês", "Semantic", 1),
    @(99, 98, "29/12/2025 21:58", "879", $true,  "Inglês", "Interpretação de Texto", 1),
    @(100, 99, "29/12/2025 22:04", "878", $true,  "Inglês", "Interpretação de Texto", 1),
    @(101, 100, "29/12/2025 22:11", 970, $false, "Inglês", "Semantic", 0)
)

foreach ($r in $rows) {
    $rowNum     = $r[0]
    $idVal      = $r[1]
    $dataVal    = $r[2]
    $questaoVal = $r[3]
    $questaoIsText = $r[4]
    $disciplina = $r[5]
    $assunto    = $r[6]
    $resultado  = $r[7]

    $ws.Cells.Item($rowNum, 1).Value = $idVal
    $ws.Cells.Item($rowNum, 2).Value = $dataVal

    $cQuestao = $ws.Cells.Item($rowNum, 3)
    if ($questaoIsText) {
        $cQuestao.NumberFormat = "@"
        $cQuestao.Value = $questaoVal
    } else {
        $cQuestao.Value = $questaoVal
    }

    $ws.Cells.Item($rowNum, 4).Value = $disciplina
    $ws.Cells.Item($rowNum, 5).Value = $assunto
    $ws.Cells.Item($rowNum, 6).Value = $resultado
}
